$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value  = "#PBM:DWP:Request:Check user type; must be Customer or Employee.:Educate"
$ws.Range("E3").Value  = "#PBM:DWP:Data:Investigate data sync issue between Remedy and Service Desk.:R&D"
$ws.Range("E4").Value  = "#PBM:DWP:Customization:Error editing questionnaire; check customization settings.:NA"
$ws.Range("E5").Value  = "#PBM:DWPC:Permissions:Granted administrator permissions to specified users.:Customization"
$ws.Range("E6").Value  = "#PBM:DWP:Integration:Guide on automating user creation via API for DWP access.:Educate"
$ws.Range("E7").Value  = "#PBM:DWP:Defect:Date format inconsistency; consult R&D for expected behavior.:R&D"
$ws.Range("E8").Value  = "#PBM:DWP:Request:Investigate request cancellation delay issue.:R&D"
$ws.Range("E9").Value  = "#PBM:DWP:Workflow:Investigate workflow approval delay and request creation issue.:R&D"
$ws.Range("E10").Value = "#PBM:SRM:Report:Guide on mapping SRM fields to Helix Dashboard equivalents.:Educate"
$ws.Range("E11").Value = "#PBM:DWPC:Data:Provided guidance on credential encryption/decryption in DWPC.:Educate"
$ws.Range("E12").Value = "#PBM:DWP:Notification Template:Disabled 'Waiting approval' notification via DB query.:Customization"
$ws.Range("E13").Value = "#PBM:DWP:Configuration:Increased autoclose duration to 1 year post-upgrade.:Customization"
$ws.Range("E14").Value = "#PBM:SRM:Workflow:Investigate SR status sync issue with completed WOs.:R&D"
$ws.Range("E15").Value = "#PBM:DWP:Multitenancy:Setup enhanced catalog for sub tenant on production.:Customization"
$ws.Range("E16").Value = "#PBM:DWPC:Log:Investigate logs for deactivation cause.:NA"
$ws.Range("E17").Value = "#PBM:DWP:Request:Manually updated request status to completed.:NA"
$ws.Range("E18").Value = "#PBM:DWP:Broadcast:Broadcast not fully removed; sync issue suspected.:NA"
$ws.Range("E19").Value = "#PBM:DWP:Configuration:Check reCAPTCHA v3 configuration and adjust settings.:Educate"
$ws.Range("E20").Value = "#PBM:DWP:Customization:Educated on customizing comment display settings.:Educate"
$ws.Range("E21").Value = "#PBM:DWPC:Configuration:Page error when enabling services in chatbot; check configuration settings.:NA"
